$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.841.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.469.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.06%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.468.26'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.130'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.83%  '

$ws.Range("E11").Value = '  -2.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.407'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.071.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.900.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000169'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.470.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.37%  '

$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.531'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.93%  '

$ws.Range("E28").Value = '  +1.39%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.96%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.89%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.27%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '29.15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.891'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.791.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.98%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.25%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.65%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.62%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0680'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '326.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0286'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.807'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.79%  '
